# Reference-doc style refresh: Title/TitleChar lose their condensed
# character spacing + kerning threshold, and Author/Date become
# Title-based styles with their own (smaller) run size instead of a
# style-local center alignment override.

$d = $word.ActiveDocument

# --- Title / TitleChar: drop the -10 twip char spacing and 28 half-pt
#     kerning threshold from the run properties (both styles carried
#     identical rPr tweaks). ---
foreach ($styleName in @("Title", "TitleChar")) {
    $s = $d.Styles($styleName)
    $s.Font.Spacing = 0
    $s.Font.Kerning = 0
}

# --- Author / Date: base them on Title (inherits the centered
#     paragraph alignment instead of a local override) and give them
#     their own 12pt run size. ---
$titleStyle = $d.Styles("Title")
foreach ($styleName in @("Author", "Date")) {
    $s = $d.Styles($styleName)
    $s.BaseStyle = $titleStyle
    $s.Font.Size = 12
    $s.Font.SizeBi = 12
}
